$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value = 11311.875
$ws.Range("I6").Value = 12062.733
$ws.Range("K6").Value = 36188.199
$ws.Range("M6").Value = -36076.199
$ws.Range("H28").Value = 279
$ws.Range("I28").Value = 238.18182
$ws.Range("J28").Value = 503.5
$ws.Range("K28").Value = 238.18182
$ws.Range("L28").Value = 503.5
$ws.Range("M28").Value = 246.81818
$ws.Range("N28").Value = -1473.5
$ws.Range("H33").Value = 2373943.8
$ws.Range("J33").Value = 350
$ws.Range("L33").Value = 350
$ws.Range("N33").Value = -808
$ws.Range("H57").Value = 49998
$ws.Range("J57").Value = 49998
$ws.Range("L57").Value = 149994
$ws.Range("N57").Value = -150992
$ws.Range("H62").Value = 1200
$ws.Range("I62").Value = 1200
$ws.Range("J62").Value = 1200
$ws.Range("K62").Value = 1200
$ws.Range("L62").Value = 1200
$ws.Range("M62").Value = -576
$ws.Range("N62").Value = -2448
$ws.Range("H65").Value = 1200
$ws.Range("I65").Value = 1200
$ws.Range("J65").Value = 1200
$ws.Range("K65").Value = 6000
$ws.Range("L65").Value = 6000
$ws.Range("M65").Value = -2880
$ws.Range("N65").Value = -12240
$ws.Range("H80").Value = 981.087
$ws.Range("I80").Value = 534.2222
$ws.Range("J80").Value = 1268.3572
$ws.Range("K80").Value = 1602.6666
$ws.Range("L80").Value = 3805.0716
$ws.Range("M80").Value = -604.6666
$ws.Range("N80").Value = -5801.071599999999
$ws.Range("H83").Value = 981.087
$ws.Range("I83").Value = 534.2222
$ws.Range("J83").Value = 1268.3572
$ws.Range("K83").Value = 4807.999800000001
$ws.Range("L83").Value = 11415.2148
$ws.Range("M83").Value = 184.0001999999995
$ws.Range("N83").Value = -21399.2148
$ws.Range("H87").Value = 50011.5
$ws.Range("J87").Value = 50011.5
$ws.Range("L87").Value = 50011.5
$ws.Range("N87").Value = -52507.5
$ws.Range("H90").Value = 50011.5
$ws.Range("J90").Value = 50011.5
$ws.Range("L90").Value = 150034.5
$ws.Range("N90").Value = -162514.5
$ws.Range("H94").Value = 1127.2
$ws.Range("I94").Value = 1127.2
$ws.Range("K94").Value = 1127.2
$ws.Range("M94").Value = -676.2
$ws.Range("H98").Value = 1109.4166
$ws.Range("I98").Value = 1047.8636
$ws.Range("K98").Value = 1047.8636
$ws.Range("M98").Value = 450.1364000000001
$ws.Range("H99").Value = 12871.375
$ws.Range("I99").Value = 328.5
$ws.Range("J99").Value = 50500
$ws.Range("K99").Value = 985.5
$ws.Range("L99").Value = 151500
$ws.Range("M99").Value = 512.5
$ws.Range("N99").Value = -154496
$ws.Range("H100").Value = 92271.09
$ws.Range("I100").Value = 126365.75
$ws.Range("J100").Value = 1352
$ws.Range("K100").Value = 126365.75
$ws.Range("L100").Value = 1352
$ws.Range("M100").Value = -125824.75
$ws.Range("N100").Value = -2434
$ws.Range("H101").Value = 2765.6667
$ws.Range("I101").Value = 303.1
$ws.Range("J101").Value = 7690.8
$ws.Range("K101").Value = 909.3000000000001
$ws.Range("L101").Value = 23072.4
$ws.Range("M101").Value = 712.6999999999999
$ws.Range("N101").Value = -26316.4
$ws.Range("H107").Value = 1657.9615
$ws.Range("I107").Value = 1380
$ws.Range("K107").Value = 1380
$ws.Range("M107").Value = 540
$ws.Range("H113").Value = 6050.625
$ws.Range("I113").Value = 6652.25
$ws.Range("K113").Value = 6652.25
$ws.Range("M113").Value = -3398.25
$ws.Range("H115").Value = 424.55554
$ws.Range("J115").Value = 488.5
$ws.Range("L115").Value = 1465.5
$ws.Range("N115").Value = -4599.5
$ws.Range("H122").Value = 1109.4166
$ws.Range("I122").Value = 1047.8636
$ws.Range("K122").Value = 3143.5908
$ws.Range("M122").Value = -693.5907999999999
$ws.Range("H132").Value = 3574184.5
$ws.Range("I132").Value = 4653525.5
$ws.Range("J132").Value = 4055.6155
$ws.Range("K132").Value = 13960576.5
$ws.Range("L132").Value = 12166.8465
$ws.Range("M132").Value = -13958046.5
$ws.Range("N132").Value = -17226.8465
$ws.Range("H137").Value = 8345.758
$ws.Range("I137").Value = 4594.7144
$ws.Range("J137").Value = 14910.083
$ws.Range("K137").Value = 13784.1432
$ws.Range("L137").Value = 44730.249
$ws.Range("M137").Value = -11234.1432
$ws.Range("N137").Value = -49830.249
$ws.Range("H138").Value = 4600.729
$ws.Range("I138").Value = 1960.1111
$ws.Range("J138").Value = 5760.0244
$ws.Range("K138").Value = 5880.3333
$ws.Range("L138").Value = 17280.0732
$ws.Range("M138").Value = -740.3333000000002
$ws.Range("N138").Value = -27560.0732

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H17").Value = 2766.7693
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2766.7693
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = 2766.7693
$ws.Range("N17").Value = -3112.7693
$ws.Range("H29").Value = 18828.834
$ws.Range("J29").Value = 18995.295
$ws.Range("L29").Value = 18995.295
$ws.Range("N29").Value = -19611.295
$ws.Range("H32").Value = 2339.6182
$ws.Range("I32").Value = 1458.551
$ws.Range("K32").Value = 1458.551
$ws.Range("M32").Value = -1171.551
$ws.Range("H45").Value = 9486.071
$ws.Range("I45").Value = 13236.667
$ws.Range("K45").Value = 13236.667
$ws.Range("M45").Value = -12859.667
$ws.Range("H60").Value = 40768.152
$ws.Range("J60").Value = 40768.152
$ws.Range("L60").Value = 40768.152
$ws.Range("N60").Value = -42234.152
$ws.Range("H61").Value = 4766.6
$ws.Range("I61").Value = 2670.258
$ws.Range("J61").Value = 11987.333
$ws.Range("K61").Value = 2670.258
$ws.Range("L61").Value = 11987.333
$ws.Range("M61").Value = -2458.258
$ws.Range("N61").Value = -12411.333
$ws.Range("H63").Value = 2138.1428
$ws.Range("I63").Value = 1244.5
$ws.Range("K63").Value = 1244.5
$ws.Range("M63").Value = -558.5
$ws.Range("H66").Value = 2138.1428
$ws.Range("I66").Value = 1244.5
$ws.Range("K66").Value = 6222.5
$ws.Range("M66").Value = -2790.5
$ws.Range("H74").Value = 12348952
$ws.Range("I74").Value = 17544990
$ws.Range("J74").Value = 8360.875
$ws.Range("K74").Value = 17544990
$ws.Range("L74").Value = 8360.875
$ws.Range("M74").Value = -17544116
$ws.Range("N74").Value = -10108.875
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("H76").Value = 17317.5
$ws.Range("J76").Value = 17317.5
$ws.Range("L76").Value = 17317.5
$ws.Range("N76").Value = -17993.5
$ws.Range("H77").Value = 12348952
$ws.Range("I77").Value = 17544990
$ws.Range("J77").Value = 8360.875
$ws.Range("K77").Value = 87724950
$ws.Range("L77").Value = 41804.375
$ws.Range("M77").Value = -87720582
$ws.Range("N77").Value = -50540.375
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("H79").Value = 17317.5
$ws.Range("J79").Value = 17317.5
$ws.Range("L79").Value = 17317.5
$ws.Range("N79").Value = -19657.5
$ws.Range("H96").Value = 83999.5
$ws.Range("J96").Value = 83999.5
$ws.Range("L96").Value = 83999.5
$ws.Range("N96").Value = -89491.5
$ws.Range("H97").Value = 1869.8
$ws.Range("I97").Value = 2222.75
$ws.Range("J97").Value = 458
$ws.Range("K97").Value = 2222.75
$ws.Range("L97").Value = 458
$ws.Range("M97").Value = -1726.75
$ws.Range("N97").Value = -1450
$ws.Range("H102").Value = 3019.6667
$ws.Range("I102").Value = 3029.6155
$ws.Range("K102").Value = 3029.6155
$ws.Range("M102").Value = -1407.6155
$ws.Range("H104").Value = 84498.5
$ws.Range("J104").Value = 84498.5
$ws.Range("L104").Value = 84498.5
$ws.Range("N104").Value = -91486.5
$ws.Range("H111").Value = 116533
$ws.Range("J111").Value = 116533
$ws.Range("L111").Value = 116533
$ws.Range("N111").Value = -124713
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").Value = 0
$ws.Range("H122").Value = 2695
$ws.Range("I122").Value = 2695
$ws.Range("K122").Value = 8085
$ws.Range("M122").Value = -5635
$ws.Range("H132").Value = 12765.944
$ws.Range("I132").Value = 12809.381
$ws.Range("J132").Value = 12705.134
$ws.Range("K132").Value = 38428.143
$ws.Range("L132").Value = 38115.402
$ws.Range("M132").Value = -35898.143
$ws.Range("N132").Value = -43175.402
$ws.Range("H136").Value = 4766.6
$ws.Range("I136").Value = 2670.258
$ws.Range("J136").Value = 11987.333
$ws.Range("K136").Value = 8010.773999999999
$ws.Range("L136").Value = 35961.999
$ws.Range("M136").Value = -5460.773999999999
$ws.Range("N136").Value = -41061.999
$ws.Range("H139").Value = 344999.5
$ws.Range("J139").Value = 344999.5
$ws.Range("L139").Value = 344999.5
$ws.Range("N139").Value = -355279.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H74").Value = 59927.5
$ws.Range("J74").Value = 59927.5
$ws.Range("L74").Value = 59927.5
$ws.Range("N74").Value = -61799.5
$ws.Range("H77").Value = 59927.5
$ws.Range("J77").Value = 59927.5
$ws.Range("L77").Value = 179782.5
$ws.Range("N77").Value = -189142.5
$ws.Range("H92").Value = 40100
$ws.Range("J92").Value = 40100
$ws.Range("L92").Value = 40100
$ws.Range("N92").Value = -45092
$ws.Range("H99").Value = 3717.2917
$ws.Range("I99").Value = 3773.4092
$ws.Range("K99").Value = 3773.4092
$ws.Range("M99").Value = -2275.4092
$ws.Range("H105").Value = 3863.647
$ws.Range("J105").Value = 2791.2856
$ws.Range("L105").Value = 2791.2856
$ws.Range("N105").Value = -6285.2856
$ws.Range("H114").Value = 76660.336
$ws.Range("J114").Value = 76660.336
$ws.Range("L114").Value = 76660.336
$ws.Range("N114").Value = -85338.336
$ws.Range("H132").Value = 79759.75
$ws.Range("J132").Value = 79759.75
$ws.Range("L132").Value = 79759.75
$ws.Range("N132").Value = -89879.75
$ws.Range("H134").Value = 7944.2666
$ws.Range("I134").Value = 4653.524
$ws.Range("J134").Value = 15622.667
$ws.Range("K134").Value = 13960.572
$ws.Range("L134").Value = 46868.001
$ws.Range("M134").Value = -11425.572
$ws.Range("N134").Value = -51938.001

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 172.05882
$ws.Range("I7").Value = 210.07692
$ws.Range("K7").Value = 210.07692
$ws.Range("M7").Value = -97.07692
$ws.Range("H15").Value = 12278.95
$ws.Range("J15").Value = 14945.5625
$ws.Range("L15").Value = 14945.5625
$ws.Range("N15").Value = -15285.5625
$ws.Range("H22").Value = 271.45834
$ws.Range("I22").Value = 282.1
$ws.Range("J22").Value = 218.25
$ws.Range("K22").Value = 282.1
$ws.Range("L22").Value = 218.25
$ws.Range("M22").Value = 67.89999999999998
$ws.Range("N22").Value = -918.25
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("H62").Value = 6298.6
$ws.Range("I62").Value = 5252
$ws.Range("J62").Value = 6996.3335
$ws.Range("K62").Value = 5252
$ws.Range("L62").Value = 6996.3335
$ws.Range("M62").Value = -4628
$ws.Range("N62").Value = -8244.333500000001
$ws.Range("H65").Value = 6298.6
$ws.Range("I65").Value = 5252
$ws.Range("J65").Value = 6996.3335
$ws.Range("K65").Value = 26260
$ws.Range("L65").Value = 34981.6675
$ws.Range("M65").Value = -23140
$ws.Range("N65").Value = -41221.6675
$ws.Range("H86").Value = 4636.222
$ws.Range("I86").Value = 4746.5713
$ws.Range("K86").Value = 4746.5713
$ws.Range("M86").Value = -3623.5713
$ws.Range("H89").Value = 4636.222
$ws.Range("I89").Value = 4746.5713
$ws.Range("K89").Value = 23732.8565
$ws.Range("M89").Value = -18116.8565
$ws.Range("H107").Value = 34170.332
$ws.Range("J107").Value = 50555.5
$ws.Range("L107").Value = 50555.5
$ws.Range("N107").Value = -54395.5
$ws.Range("H132").Value = 17251.559
$ws.Range("I132").Value = 10836.76
$ws.Range("K132").Value = 32510.28
$ws.Range("M132").Value = -29980.28

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H3").Value = 3727.5
$ws.Range("I3").Value = 3727.5
$ws.Range("K3").Value = 11182.5
$ws.Range("M3").Value = -11070.5
$ws.Range("H4").Value = 81041070
$ws.Range("I4").Value = 81041070
$ws.Range("K4").Value = 243123210
$ws.Range("M4").Value = -243123098
$ws.Range("I7").Value = 62500136
$ws.Range("K7").Value = 187500408
$ws.Range("M7").Value = -187500296
$ws.Range("H8").Value = 585.1429000000001
$ws.Range("I8").Value = 585.1429000000001
$ws.Range("K8").Value = 1755.4287
$ws.Range("M8").Value = -1616.4287
$ws.Range("H26").Value = 78
$ws.Range("I26").Value = 53.6
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 160.8
$ws.Range("L26").Value = 600
$ws.Range("M26").Value = 127.2
$ws.Range("N26").Value = -1176
$ws.Range("H34").Value = 2229.3794
$ws.Range("J34").Value = 2884
$ws.Range("L34").Value = 8652
$ws.Range("N34").Value = -8820
$ws.Range("H45").Value = 8375
$ws.Range("I45").Value = 7750
$ws.Range("K45").Value = 23250
$ws.Range("M45").Value = -22718
$ws.Range("H68").Value = 1340.5834
$ws.Range("J68").Value = 1425.091
$ws.Range("L68").Value = 4275.272999999999
$ws.Range("N68").Value = -5897.272999999999
$ws.Range("H71").Value = 1340.5834
$ws.Range("J71").Value = 1425.091
$ws.Range("L71").Value = 12825.819
$ws.Range("N71").Value = -20937.819
$ws.Range("H75").Value = 897.7143
$ws.Range("I75").Value = 766
$ws.Range("J75").Value = 1073.3334
$ws.Range("K75").Value = 2298
$ws.Range("L75").Value = 3220.0002
$ws.Range("M75").Value = -1300
$ws.Range("N75").Value = -5216.0002
$ws.Range("H78").Value = 897.7143
$ws.Range("I78").Value = 766
$ws.Range("J78").Value = 1073.3334
$ws.Range("K78").Value = 6894
$ws.Range("L78").Value = 9660.000599999999
$ws.Range("M78").Value = -1902
$ws.Range("N78").Value = -19644.0006
$ws.Range("H80").Value = 21987
$ws.Range("I80").Value = 100000
$ws.Range("J80").Value = 10842.286
$ws.Range("K80").Value = 300000
$ws.Range("L80").Value = 32526.858
$ws.Range("M80").Value = -299064
$ws.Range("N80").Value = -34398.858
$ws.Range("H83").Value = 21987
$ws.Range("I83").Value = 100000
$ws.Range("J83").Value = 10842.286
$ws.Range("K83").Value = 900000
$ws.Range("L83").Value = 97580.57399999999
$ws.Range("M83").Value = -895320
$ws.Range("N83").Value = -106940.574
$ws.Range("H86").Value = 817.05554
$ws.Range("I86").Value = 780.8889
$ws.Range("J86").Value = 853.2222
$ws.Range("K86").Value = 2342.6667
$ws.Range("L86").Value = 2559.6666
$ws.Range("M86").Value = -1156.6667
$ws.Range("N86").Value = -4931.6666
$ws.Range("H89").Value = 817.05554
$ws.Range("I89").Value = 780.8889
$ws.Range("J89").Value = 853.2222
$ws.Range("K89").Value = 7028.0001
$ws.Range("L89").Value = 7678.999800000001
$ws.Range("M89").Value = -1100.0001
$ws.Range("N89").Value = -19534.9998
$ws.Range("H92").Value = 466.33334
$ws.Range("I92").Value = 99
$ws.Range("K92").Value = 297
$ws.Range("M92").Value = 951
$ws.Range("H99").Value = 8351.200000000001
$ws.Range("J99").Value = 12199.4
$ws.Range("L99").Value = 36598.2
$ws.Range("N99").Value = -41090.2
$ws.Range("H107").Value = 955.6279
$ws.Range("I107").Value = 425.25
$ws.Range("J107").Value = 1010.02563
$ws.Range("K107").Value = 1275.75
$ws.Range("L107").Value = 3030.07689
$ws.Range("M107").Value = 644.25
$ws.Range("N107").Value = -6870.07689
$ws.Range("H108").Value = 1271.7778
$ws.Range("I108").Value = 1271.7778
$ws.Range("K108").Value = 3815.3334
$ws.Range("M108").Value = -935.3334000000004
$ws.Range("H114").Value = 507.57144
$ws.Range("I114").Value = 521.3333
$ws.Range("J114").Value = 425
$ws.Range("K114").Value = 1563.9999
$ws.Range("L114").Value = 1275
$ws.Range("M114").Value = 1690.0001
$ws.Range("N114").Value = -7783
$ws.Range("H122").Value = 1996
$ws.Range("I122").Value = 1995.5
$ws.Range("K122").Value = 17959.5
$ws.Range("M122").Value = -15509.5
$ws.Range("H126").Value = 4990.3335
$ws.Range("I126").Value = 4990.3335
$ws.Range("K126").Value = 14971.0005
$ws.Range("M126").Value = -10031.0005
$ws.Range("H127").Value = 4266.857
$ws.Range("J127").Value = 4266.857
$ws.Range("L127").Value = 12800.571
$ws.Range("N127").Value = -22720.571
$ws.Range("H132").Value = 1828.6666
$ws.Range("I132").Value = 1350.75
$ws.Range("J132").Value = 1965.2142
$ws.Range("K132").Value = 12156.75
$ws.Range("L132").Value = 17686.9278
$ws.Range("M132").Value = -9626.75
$ws.Range("N132").Value = -22746.9278
$ws.Range("H139").Value = 487.07144
$ws.Range("I139").Value = 487.07144
$ws.Range("K139").Value = 1461.21432
$ws.Range("M139").Value = 3678.78568
$ws.Range("H140").Value = 1065.7742
$ws.Range("I140").Value = 1001.3
$ws.Range("K140").Value = 3003.9
$ws.Range("M140").Value = 2176.1
$ws.Range("H141").Value = 3431.6667
$ws.Range("I141").Value = 1860.625
$ws.Range("K141").Value = 5581.875
$ws.Range("M141").Value = -401.875

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 93.5625
$ws.Range("I2").Value = 67.13333
$ws.Range("J2").Value = 490
$ws.Range("K2").Value = 67.13333
$ws.Range("L2").Value = 490
$ws.Range("M2").Value = 45.86667
$ws.Range("N2").Value = -716
$ws.Range("H18").Value = 14998.333
$ws.Range("I18").Value = 14990
$ws.Range("K18").Value = 14990
$ws.Range("M18").Value = -14697
$ws.Range("H29").Value = 3500
$ws.Range("J29").Value = 3500
$ws.Range("L29").Value = 3500
$ws.Range("N29").Value = -4080
$ws.Range("H57").Value = 38683
$ws.Range("J57").Value = 38585
$ws.Range("L57").Value = 38585
$ws.Range("N57").Value = -40225
$ws.Range("H80").Value = 8226
$ws.Range("I80").Value = 1717.8
$ws.Range("K80").Value = 1717.8
$ws.Range("M80").Value = -719.8
$ws.Range("H83").Value = 8226
$ws.Range("I83").Value = 1717.8
$ws.Range("K83").Value = 8589
$ws.Range("M83").Value = -3597
$ws.Range("H97").Value = 2781.9167
$ws.Range("I97").Value = 2636.125
$ws.Range("K97").Value = 2636.125
$ws.Range("M97").Value = -2140.125
$ws.Range("H102").Value = 2757.04
$ws.Range("I102").Value = 2780.261
$ws.Range("K102").Value = 2780.261
$ws.Range("M102").Value = -1158.261
$ws.Range("H132").Value = 13440.167
$ws.Range("I132").Value = 7577.636
$ws.Range("J132").Value = 22652.715
$ws.Range("K132").Value = 22732.908
$ws.Range("L132").Value = 67958.145
$ws.Range("M132").Value = -20202.908
$ws.Range("N132").Value = -73018.145

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H20").Value = 41493.617
$ws.Range("I20").Value = 999
$ws.Range("J20").Value = 44868.168
$ws.Range("K20").Value = 999
$ws.Range("L20").Value = 44868.168
$ws.Range("M20").Value = -773
$ws.Range("N20").Value = -45320.168
$ws.Range("H22").Value = 1434.3334
$ws.Range("I22").Value = 1647.6364
$ws.Range("J22").Value = 847.75
$ws.Range("K22").Value = 1647.6364
$ws.Range("L22").Value = 847.75
$ws.Range("M22").Value = -1352.6364
$ws.Range("N22").Value = -1437.75
$ws.Range("H23").Value = 22417.666
$ws.Range("I23").Value = 26601.2
$ws.Range("J23").Value = 1500
$ws.Range("K23").Value = 26601.2
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = -26371.2
$ws.Range("N23").Value = -1960
$ws.Range("H25").Value = 10002
$ws.Range("I25").Value = 10002
$ws.Range("K25").Value = 10002
$ws.Range("M25").Value = -9772
$ws.Range("H27").Value = 1434.3334
$ws.Range("I27").Value = 1647.6364
$ws.Range("J27").Value = 847.75
$ws.Range("K27").Value = 1647.6364
$ws.Range("L27").Value = 847.75
$ws.Range("M27").Value = -1540.6364
$ws.Range("N27").Value = -1061.75
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10912
$ws.Range("H46").Value = 1268.2444
$ws.Range("I46").Value = 986.9545000000001
$ws.Range("K46").Value = 986.9545000000001
$ws.Range("M46").Value = -798.9545000000001
$ws.Range("H61").Value = 3905.2856
$ws.Range("I61").Value = 3153.4546
$ws.Range("J61").Value = 6662
$ws.Range("K61").Value = 3153.4546
$ws.Range("L61").Value = 6662
$ws.Range("M61").Value = -2951.4546
$ws.Range("N61").Value = -7066
$ws.Range("H68").Value = 3233.8125
$ws.Range("I68").Value = 1235.1
$ws.Range("K68").Value = 1235.1
$ws.Range("M68").Value = -486.0999999999999
$ws.Range("H71").Value = 3233.8125
$ws.Range("I71").Value = 1235.1
$ws.Range("K71").Value = 6175.5
$ws.Range("M71").Value = -2431.5
$ws.Range("H82").Value = 1277.875
$ws.Range("I82").Value = 966.2
$ws.Range("K82").Value = 966.2
$ws.Range("M82").Value = -605.2
$ws.Range("H85").Value = 1277.875
$ws.Range("I85").Value = 966.2
$ws.Range("K85").Value = 966.2
$ws.Range("M85").Value = 281.8
$ws.Range("H88").Value = 32749.75
$ws.Range("I88").Value = 15500
$ws.Range("J88").Value = 49999.5
$ws.Range("K88").Value = 15500
$ws.Range("L88").Value = 49999.5
$ws.Range("M88").Value = -15072
$ws.Range("N88").Value = -50855.5
$ws.Range("H91").Value = 32749.75
$ws.Range("I91").Value = 15500
$ws.Range("J91").Value = 49999.5
$ws.Range("K91").Value = 15500
$ws.Range("L91").Value = 49999.5
$ws.Range("M91").Value = -14018
$ws.Range("N91").Value = -52963.5
$ws.Range("H93").Value = 4101.5356
$ws.Range("I93").Value = 4237.2856
$ws.Range("J93").Value = 3694.2856
$ws.Range("K93").Value = 4237.2856
$ws.Range("L93").Value = 3694.2856
$ws.Range("M93").Value = -2989.2856
$ws.Range("N93").Value = -6190.2856
$ws.Range("H113").Value = 3905.2856
$ws.Range("I113").Value = 3153.4546
$ws.Range("J113").Value = 6662
$ws.Range("K113").Value = 3153.4546
$ws.Range("L113").Value = 6662
$ws.Range("M113").Value = -983.4546
$ws.Range("N113").Value = -11002
$ws.Range("H136").Value = 4117356.8
$ws.Range("I136").Value = 4275674
$ws.Range("K136").Value = 12827022
$ws.Range("M136").Value = -12824472

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 2352.7727
$ws.Range("J4").Value = 2092.8948
$ws.Range("L4").Value = 2092.8948
$ws.Range("N4").Value = -2318.8948
$ws.Range("H14").Value = 10874.5
$ws.Range("J14").Value = 9570.857
$ws.Range("L14").Value = 9570.857
$ws.Range("N14").Value = -9906.857
$ws.Range("H29").Value = 16005
$ws.Range("J29").Value = 12011
$ws.Range("L29").Value = 12011
$ws.Range("N29").Value = -12591
$ws.Range("H58").Value = 8000
$ws.Range("I58").Value = 8000
$ws.Range("K58").Value = 8000
$ws.Range("M58").Value = -7692
$ws.Range("H62").Value = 3200
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 3200
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H81").Value = 2632.182
$ws.Range("I81").Value = 2325.6667
$ws.Range("K81").Value = 4651.3334
$ws.Range("M81").Value = -3590.3334
$ws.Range("H84").Value = 2632.182
$ws.Range("I84").Value = 2325.6667
$ws.Range("K84").Value = 23256.667
$ws.Range("M84").Value = -17952.667
$ws.Range("H93").Value = 56794.332
$ws.Range("J93").Value = 56794.332
$ws.Range("L93").Value = 56794.332
$ws.Range("N93").Value = -61786.332
$ws.Range("H95").Value = 25166.666
$ws.Range("J95").Value = 25166.666
$ws.Range("L95").Value = 25166.666
$ws.Range("N95").Value = -30658.666
$ws.Range("H99").Value = 49563.332
$ws.Range("J99").Value = 49563.332
$ws.Range("L99").Value = 49563.332
$ws.Range("N99").Value = -55553.332
$ws.Range("H119").Value = 77800
$ws.Range("J119").Value = 77800
$ws.Range("L119").Value = 77800
$ws.Range("N119").Value = -87476
$ws.Range("H121").Value = 95000
$ws.Range("J121").Value = 95000
$ws.Range("L121").Value = 95000
$ws.Range("N121").Value = -98494
$ws.Range("H122").Value = 7530.5557
$ws.Range("J122").Value = 3500
$ws.Range("L122").Value = 10500
$ws.Range("N122").Value = -15400
$ws.Range("H136").Value = 2236.0981
$ws.Range("I136").Value = 1973.7046
$ws.Range("K136").Value = 5921.1138
$ws.Range("M136").Value = -3371.1138
